# Update "want to go" counts (column F) on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes   = $wb.Worksheets.Item(4)   # 全部类型

# Sheet "展览" (sheet1)
$wsExhibition.Range("F2").Value = 5549
$wsExhibition.Range("F3").Value = 625
$wsExhibition.Range("F4").Value = 12556
$wsExhibition.Range("F6").Value = 620
$wsExhibition.Range("F7").Value = 192
$wsExhibition.Range("F8").Value = 369
$wsExhibition.Range("F9").Value = 1162

# Sheet "全部类型" (sheet4)
$wsAllTypes.Range("F2").Value = 5549
$wsAllTypes.Range("F3").Value = 625
$wsAllTypes.Range("F5").Value = 12556
$wsAllTypes.Range("F7").Value = 620
$wsAllTypes.Range("F8").Value = 192
$wsAllTypes.Range("F11").Value = 369
$wsAllTypes.Range("F12").Value = 1162
